$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(281, 4).Value = 44636
$ws.Cells.Item(281, 12).Value = "Primera"
$ws.Cells.Item(281, 13).Value = 200
$ws.Cells.Item(281, 14).Value = 4400
$ws.Cells.Item(281, 15).Value = 4500
$ws.Cells.Item(281, 16).Value = 4450
$ws.Cells.Item(281, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(281, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(281, 19).Value = 4450
$ws.Cells.Item(281, 20).Value = 1
$ws.Cells.Item(282, 4).Value = 44636
$ws.Cells.Item(282, 12).Value = "Segunda"
$ws.Cells.Item(282, 13).Value = 100
$ws.Cells.Item(282, 14).Value = 3800
$ws.Cells.Item(282, 15).Value = 3800
$ws.Cells.Item(282, 16).Value = 3800
$ws.Cells.Item(282, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(282, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(282, 19).Value = 3800
$ws.Cells.Item(282, 20).Value = 1
$ws.Cells.Item(283, 4).Value = 44411
$ws.Cells.Item(283, 12).Value = "Especial"
$ws.Cells.Item(283, 13).Value = 150
$ws.Cells.Item(283, 14).Value = 35000
$ws.Cells.Item(283, 15).Value = 35000
$ws.Cells.Item(283, 16).Value = 35000
$ws.Cells.Item(283, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(283, 18).Value = "Perú"
$ws.Cells.Item(283, 19).Value = 3500
$ws.Cells.Item(283, 20).Value = 10
$ws.Cells.Item(284, 4).Value = 44411
$ws.Cells.Item(284, 12).Value = "Primera"
$ws.Cells.Item(284, 13).Value = 200
$ws.Cells.Item(284, 14).Value = 30000
$ws.Cells.Item(284, 15).Value = 30000
$ws.Cells.Item(284, 16).Value = 30000
$ws.Cells.Item(284, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(284, 18).Value = "Perú"
$ws.Cells.Item(284, 19).Value = 3000
$ws.Cells.Item(284, 20).Value = 10
$ws.Cells.Item(285, 4).Value = 44281
$ws.Cells.Item(285, 12).Value = "Primera"
$ws.Cells.Item(285, 13).Value = 50
$ws.Cells.Item(285, 14).Value = 6400
$ws.Cells.Item(285, 15).Value = 6400
$ws.Cells.Item(285, 16).Value = 6400
$ws.Cells.Item(285, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(285, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(285, 19).Value = 6400
$ws.Cells.Item(285, 20).Value = 1
$ws.Cells.Item(286, 4).Value = 44271
$ws.Cells.Item(286, 12).Value = "Primera"
$ws.Cells.Item(286, 13).Value = 200
$ws.Cells.Item(286, 14).Value = 5800
$ws.Cells.Item(286, 15).Value = 6000
$ws.Cells.Item(286, 16).Value = 5900
$ws.Cells.Item(286, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(286, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(286, 19).Value = 5900
$ws.Cells.Item(286, 20).Value = 1
$ws.Cells.Item(287, 4).Value = 44271
$ws.Cells.Item(287, 12).Value = "Segunda"
$ws.Cells.Item(287, 13).Value = 100
$ws.Cells.Item(287, 14).Value = 4900
$ws.Cells.Item(287, 15).Value = 4900
$ws.Cells.Item(287, 16).Value = 4900
$ws.Cells.Item(287, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(287, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(287, 19).Value = 4900
$ws.Cells.Item(287, 20).Value = 1
$ws.Cells.Item(288, 4).Value = 44554
$ws.Cells.Item(288, 12).Value = "Primera"
$ws.Cells.Item(288, 13).Value = 200
$ws.Cells.Item(288, 14).Value = 3900
$ws.Cells.Item(288, 15).Value = 4000
$ws.Cells.Item(288, 16).Value = 3950
$ws.Cells.Item(288, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(288, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(288, 19).Value = 3950
$ws.Cells.Item(288, 20).Value = 1
$ws.Cells.Item(289, 4).Value = 44554
$ws.Cells.Item(289, 12).Value = "Segunda"
$ws.Cells.Item(289, 13).Value = 100
$ws.Cells.Item(289, 14).Value = 3500
$ws.Cells.Item(289, 15).Value = 3500
$ws.Cells.Item(289, 16).Value = 3500
$ws.Cells.Item(289, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(289, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(289, 19).Value = 3500
$ws.Cells.Item(289, 20).Value = 1
$ws.Cells.Item(290, 4).Value = 44162
$ws.Cells.Item(290, 12).Value = "Primera"
$ws.Cells.Item(290, 13).Value = 150
$ws.Cells.Item(290, 14).Value = 4200
$ws.Cells.Item(290, 15).Value = 4200
$ws.Cells.Item(290, 16).Value = 4200
$ws.Cells.Item(290, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(290, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(290, 19).Value = 4200
$ws.Cells.Item(290, 20).Value = 1
$ws.Cells.Item(291, 4).Value = 44162
$ws.Cells.Item(291, 12).Value = "Segunda"
$ws.Cells.Item(291, 13).Value = 150
$ws.Cells.Item(291, 14).Value = 3600
$ws.Cells.Item(291, 15).Value = 3600
$ws.Cells.Item(291, 16).Value = 3600
$ws.Cells.Item(291, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(291, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(291, 19).Value = 3600
$ws.Cells.Item(291, 20).Value = 1
$ws.Cells.Item(292, 4).Value = 44162
$ws.Cells.Item(292, 12).Value = "Tercera"
$ws.Cells.Item(292, 13).Value = 150
$ws.Cells.Item(292, 14).Value = 3000
$ws.Cells.Item(292, 15).Value = 3000
$ws.Cells.Item(292, 16).Value = 3000
$ws.Cells.Item(292, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(292, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(292, 19).Value = 3000
$ws.Cells.Item(292, 20).Value = 1
$ws.Cells.Item(293, 4).Value = 44516
$ws.Cells.Item(293, 12).Value = "Primera"
$ws.Cells.Item(293, 13).Value = 400
$ws.Cells.Item(293, 14).Value = 4000
$ws.Cells.Item(293, 15).Value = 4100
$ws.Cells.Item(293, 16).Value = 4050
$ws.Cells.Item(293, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(293, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(293, 19).Value = 4050
$ws.Cells.Item(293, 20).Value = 1
$ws.Cells.Item(294, 4).Value = 44516
$ws.Cells.Item(294, 12).Value = "Segunda"
$ws.Cells.Item(294, 13).Value = 200
$ws.Cells.Item(294, 14).Value = 3500
$ws.Cells.Item(294, 15).Value = 3500
$ws.Cells.Item(294, 16).Value = 3500
$ws.Cells.Item(294, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(294, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(294, 19).Value = 3500
$ws.Cells.Item(294, 20).Value = 1
$ws.Cells.Item(295, 4).Value = 44568
$ws.Cells.Item(295, 12).Value = "Primera"
$ws.Cells.Item(295, 13).Value = 300
$ws.Cells.Item(295, 14).Value = 4000
$ws.Cells.Item(295, 15).Value = 4100
$ws.Cells.Item(295, 16).Value = 4050
$ws.Cells.Item(295, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(295, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(295, 19).Value = 4050
$ws.Cells.Item(295, 20).Value = 1
$ws.Cells.Item(296, 4).Value = 44568
$ws.Cells.Item(296, 12).Value = "Segunda"
$ws.Cells.Item(296, 13).Value = 150
$ws.Cells.Item(296, 14).Value = 3600
$ws.Cells.Item(296, 15).Value = 3600
$ws.Cells.Item(296, 16).Value = 3600
$ws.Cells.Item(296, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(296, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(296, 19).Value = 3600
$ws.Cells.Item(296, 20).Value = 1
$ws.Cells.Item(297, 4).Value = 44336
$ws.Cells.Item(297, 12).Value = "Primera"
$ws.Cells.Item(297, 13).Value = 120
$ws.Cells.Item(297, 14).Value = 7000
$ws.Cells.Item(297, 15).Value = 7200
$ws.Cells.Item(297, 16).Value = 7100
$ws.Cells.Item(297, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(297, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(297, 19).Value = 7100
$ws.Cells.Item(297, 20).Value = 1
$ws.Cells.Item(298, 4).Value = 44231
$ws.Cells.Item(298, 12).Value = "Primera"
$ws.Cells.Item(298, 13).Value = 200
$ws.Cells.Item(298, 14).Value = 5400
$ws.Cells.Item(298, 15).Value = 5500
$ws.Cells.Item(298, 16).Value = 5450
$ws.Cells.Item(298, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(298, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(298, 19).Value = 5450
$ws.Cells.Item(298, 20).Value = 1
$ws.Cells.Item(299, 4).Value = 44231
$ws.Cells.Item(299, 12).Value = "Segunda"
$ws.Cells.Item(299, 13).Value = 100
$ws.Cells.Item(299, 14).Value = 4200
$ws.Cells.Item(299, 15).Value = 4200
$ws.Cells.Item(299, 16).Value = 4200
$ws.Cells.Item(299, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(299, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(299, 19).Value = 4200
$ws.Cells.Item(299, 20).Value = 1
$ws.Cells.Item(300, 4).Value = 44565
$ws.Cells.Item(300, 12).Value = "Primera"
$ws.Cells.Item(300, 13).Value = 300
$ws.Cells.Item(300, 14).Value = 4100
$ws.Cells.Item(300, 15).Value = 4200
$ws.Cells.Item(300, 16).Value = 4150
$ws.Cells.Item(300, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(300, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(300, 19).Value = 4150
$ws.Cells.Item(300, 20).Value = 1
$ws.Cells.Item(301, 4).Value = 44565
$ws.Cells.Item(301, 12).Value = "Segunda"
$ws.Cells.Item(301, 13).Value = 100
$ws.Cells.Item(301, 14).Value = 3500
$ws.Cells.Item(301, 15).Value = 3500
$ws.Cells.Item(301, 16).Value = 3500
$ws.Cells.Item(301, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(301, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(301, 19).Value = 3500
$ws.Cells.Item(301, 20).Value = 1
$ws.Cells.Item(302, 4).Value = 44400
$ws.Cells.Item(302, 12).Value = "Especial"
$ws.Cells.Item(302, 13).Value = 100
$ws.Cells.Item(302, 14).Value = 38000
$ws.Cells.Item(302, 15).Value = 38000
$ws.Cells.Item(302, 16).Value = 38000
$ws.Cells.Item(302, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(302, 18).Value = "Perú"
$ws.Cells.Item(302, 19).Value = 3800
$ws.Cells.Item(302, 20).Value = 10
$ws.Cells.Item(303, 4).Value = 44400
$ws.Cells.Item(303, 12).Value = "Primera"
$ws.Cells.Item(303, 13).Value = 100
$ws.Cells.Item(303, 14).Value = 36000
$ws.Cells.Item(303, 15).Value = 36000
$ws.Cells.Item(303, 16).Value = 36000
$ws.Cells.Item(303, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(303, 18).Value = "Perú"
$ws.Cells.Item(303, 19).Value = 3600
$ws.Cells.Item(303, 20).Value = 10
$ws.Cells.Item(304, 4).Value = 44400
$ws.Cells.Item(304, 12).Value = "Segunda"
$ws.Cells.Item(304, 13).Value = 100
$ws.Cells.Item(304, 14).Value = 30000
$ws.Cells.Item(304, 15).Value = 30000
$ws.Cells.Item(304, 16).Value = 30000
$ws.Cells.Item(304, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(304, 18).Value = "Perú"
$ws.Cells.Item(304, 19).Value = 3000
$ws.Cells.Item(304, 20).Value = 10
$ws.Cells.Item(305, 4).Value = 44627
$ws.Cells.Item(305, 12).Value = "Primera"
$ws.Cells.Item(305, 13).Value = 200
$ws.Cells.Item(305, 14).Value = 4300
$ws.Cells.Item(305, 15).Value = 4300
$ws.Cells.Item(305, 16).Value = 4300
$ws.Cells.Item(305, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(305, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(305, 19).Value = 4300
$ws.Cells.Item(305, 20).Value = 1
$ws.Cells.Item(306, 4).Value = 44627
$ws.Cells.Item(306, 12).Value = "Segunda"
$ws.Cells.Item(306, 13).Value = 100
$ws.Cells.Item(306, 14).Value = 3900
$ws.Cells.Item(306, 15).Value = 3900
$ws.Cells.Item(306, 16).Value = 3900
$ws.Cells.Item(306, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(306, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(306, 19).Value = 3900
$ws.Cells.Item(306, 20).Value = 1
$ws.Cells.Item(307, 4).Value = 44536
$ws.Cells.Item(307, 12).Value = "Primera"
$ws.Cells.Item(307, 13).Value = 200
$ws.Cells.Item(307, 14).Value = 4000
$ws.Cells.Item(307, 15).Value = 4100
$ws.Cells.Item(307, 16).Value = 4050
$ws.Cells.Item(307, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(307, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(307, 19).Value = 4050
$ws.Cells.Item(307, 20).Value = 1
$ws.Cells.Item(308, 4).Value = 44536
$ws.Cells.Item(308, 12).Value = "Segunda"
$ws.Cells.Item(308, 13).Value = 100
$ws.Cells.Item(308, 14).Value = 3500
$ws.Cells.Item(308, 15).Value = 3500
$ws.Cells.Item(308, 16).Value = 3500
$ws.Cells.Item(308, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(308, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(308, 19).Value = 3500
$ws.Cells.Item(308, 20).Value = 1
$ws.Cells.Item(309, 4).Value = 44334
$ws.Cells.Item(309, 12).Value = "Primera"
$ws.Cells.Item(309, 13).Value = 200
$ws.Cells.Item(309, 14).Value = 44000
$ws.Cells.Item(309, 15).Value = 45000
$ws.Cells.Item(309, 16).Value = 44500
$ws.Cells.Item(309, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(309, 18).Value = "Perú"
$ws.Cells.Item(309, 19).Value = 4450
$ws.Cells.Item(309, 20).Value = 10
$ws.Cells.Item(310, 4).Value = 44334
$ws.Cells.Item(310, 12).Value = "Primera"
$ws.Cells.Item(310, 13).Value = 200
$ws.Cells.Item(310, 14).Value = 7000
$ws.Cells.Item(310, 15).Value = 7200
$ws.Cells.Item(310, 16).Value = 7100
$ws.Cells.Item(310, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(310, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(310, 19).Value = 7100
$ws.Cells.Item(310, 20).Value = 1
$ws.Cells.Item(311, 4).Value = 44193
$ws.Cells.Item(311, 12).Value = "Primera"
$ws.Cells.Item(311, 13).Value = 200
$ws.Cells.Item(311, 14).Value = 5300
$ws.Cells.Item(311, 15).Value = 5400
$ws.Cells.Item(311, 16).Value = 5350
$ws.Cells.Item(311, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(311, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(311, 19).Value = 5350
$ws.Cells.Item(311, 20).Value = 1
$ws.Cells.Item(312, 4).Value = 44193
$ws.Cells.Item(312, 12).Value = "Segunda"
$ws.Cells.Item(312, 13).Value = 100
$ws.Cells.Item(312, 14).Value = 4500
$ws.Cells.Item(312, 15).Value = 4500
$ws.Cells.Item(312, 16).Value = 4500
$ws.Cells.Item(312, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(312, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(312, 19).Value = 4500
$ws.Cells.Item(312, 20).Value = 1
$ws.Cells.Item(313, 4).Value = 44362
$ws.Cells.Item(313, 12).Value = "1a nueva(o)"
$ws.Cells.Item(313, 13).Value = 150
$ws.Cells.Item(313, 14).Value = 5400
$ws.Cells.Item(313, 15).Value = 5400
$ws.Cells.Item(313, 16).Value = 5400
$ws.Cells.Item(313, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(313, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(313, 19).Value = 5400
$ws.Cells.Item(313, 20).Value = 1
$ws.Cells.Item(314, 4).Value = 44362
$ws.Cells.Item(314, 12).Value = "2a nueva(o)"
$ws.Cells.Item(314, 13).Value = 300
$ws.Cells.Item(314, 14).Value = 4500
$ws.Cells.Item(314, 15).Value = 5300
$ws.Cells.Item(314, 16).Value = 4900
$ws.Cells.Item(314, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(314, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(314, 19).Value = 4900
$ws.Cells.Item(314, 20).Value = 1
$ws.Cells.Item(315, 4).Value = 44362
$ws.Cells.Item(315, 12).Value = "Primera"
$ws.Cells.Item(315, 13).Value = 200
$ws.Cells.Item(315, 14).Value = 42000
$ws.Cells.Item(315, 15).Value = 43000
$ws.Cells.Item(315, 16).Value = 42500
$ws.Cells.Item(315, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(315, 18).Value = "Perú"
$ws.Cells.Item(315, 19).Value = 4250
$ws.Cells.Item(315, 20).Value = 10
$ws.Cells.Item(316, 4).Value = 44365
$ws.Cells.Item(316, 12).Value = "1a nueva(o)"
$ws.Cells.Item(316, 13).Value = 200
$ws.Cells.Item(316, 14).Value = 5200
$ws.Cells.Item(316, 15).Value = 5300
$ws.Cells.Item(316, 16).Value = 5250
$ws.Cells.Item(316, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(316, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(316, 19).Value = 5250
$ws.Cells.Item(316, 20).Value = 1
$ws.Cells.Item(317, 4).Value = 44365
$ws.Cells.Item(317, 12).Value = "2a nueva(o)"
$ws.Cells.Item(317, 13).Value = 100
$ws.Cells.Item(317, 14).Value = 4500
$ws.Cells.Item(317, 15).Value = 4500
$ws.Cells.Item(317, 16).Value = 4500
$ws.Cells.Item(317, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(317, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(317, 19).Value = 4500
$ws.Cells.Item(317, 20).Value = 1
$ws.Cells.Item(318, 4).Value = 44567
$ws.Cells.Item(318, 12).Value = "Primera"
$ws.Cells.Item(318, 13).Value = 200
$ws.Cells.Item(318, 14).Value = 4100
$ws.Cells.Item(318, 15).Value = 4200
$ws.Cells.Item(318, 16).Value = 4150
$ws.Cells.Item(318, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(318, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(318, 19).Value = 4150
$ws.Cells.Item(318, 20).Value = 1
$ws.Cells.Item(319, 4).Value = 44567
$ws.Cells.Item(319, 12).Value = "Segunda"
$ws.Cells.Item(319, 13).Value = 100
$ws.Cells.Item(319, 14).Value = 3500
$ws.Cells.Item(319, 15).Value = 3500
$ws.Cells.Item(319, 16).Value = 3500
$ws.Cells.Item(319, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(319, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(319, 19).Value = 3500
$ws.Cells.Item(319, 20).Value = 1
$ws.Cells.Item(320, 4).Value = 44537
$ws.Cells.Item(320, 12).Value = "Primera"
$ws.Cells.Item(320, 13).Value = 400
$ws.Cells.Item(320, 14).Value = 4000
$ws.Cells.Item(320, 15).Value = 4100
$ws.Cells.Item(320, 16).Value = 4050
$ws.Cells.Item(320, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(320, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(320, 19).Value = 4050
$ws.Cells.Item(320, 20).Value = 1
$ws.Cells.Item(321, 4).Value = 44537
$ws.Cells.Item(321, 12).Value = "Segunda"
$ws.Cells.Item(321, 13).Value = 200
$ws.Cells.Item(321, 14).Value = 3500
$ws.Cells.Item(321, 15).Value = 3500
$ws.Cells.Item(321, 16).Value = 3500
$ws.Cells.Item(321, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(321, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(321, 19).Value = 3500
$ws.Cells.Item(321, 20).Value = 1
$ws.Cells.Item(322, 4).Value = 44553
$ws.Cells.Item(322, 12).Value = "Primera"
$ws.Cells.Item(322, 13).Value = 400
$ws.Cells.Item(322, 14).Value = 3900
$ws.Cells.Item(322, 15).Value = 4000
$ws.Cells.Item(322, 16).Value = 3950
$ws.Cells.Item(322, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(322, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(322, 19).Value = 3950
$ws.Cells.Item(322, 20).Value = 1
$ws.Cells.Item(323, 4).Value = 44553
$ws.Cells.Item(323, 12).Value = "Segunda"
$ws.Cells.Item(323, 13).Value = 200
$ws.Cells.Item(323, 14).Value = 3500
$ws.Cells.Item(323, 15).Value = 3500
$ws.Cells.Item(323, 16).Value = 3500
$ws.Cells.Item(323, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(323, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(323, 19).Value = 3500
$ws.Cells.Item(323, 20).Value = 1
$ws.Cells.Item(324, 4).Value = 44397
$ws.Cells.Item(324, 12).Value = "Especial"
$ws.Cells.Item(324, 13).Value = 200
$ws.Cells.Item(324, 14).Value = 40000
$ws.Cells.Item(324, 15).Value = 40000
$ws.Cells.Item(324, 16).Value = 40000
$ws.Cells.Item(324, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(324, 18).Value = "Perú"
$ws.Cells.Item(324, 19).Value = 4000
$ws.Cells.Item(324, 20).Value = 10
$ws.Cells.Item(325, 4).Value = 44462
$ws.Cells.Item(325, 12).Value = "1a nueva(o)"
$ws.Cells.Item(325, 13).Value = 80
$ws.Cells.Item(325, 14).Value = 4000
$ws.Cells.Item(325, 15).Value = 4000
$ws.Cells.Item(325, 16).Value = 4000
$ws.Cells.Item(325, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(325, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(325, 19).Value = 4000
$ws.Cells.Item(325, 20).Value = 1
$ws.Cells.Item(326, 4).Value = 44446
$ws.Cells.Item(326, 12).Value = "1a nueva(o)"
$ws.Cells.Item(326, 13).Value = 150
$ws.Cells.Item(326, 14).Value = 4200
$ws.Cells.Item(326, 15).Value = 4200
$ws.Cells.Item(326, 16).Value = 4200
$ws.Cells.Item(326, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(326, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(326, 19).Value = 4200
$ws.Cells.Item(326, 20).Value = 1
$ws.Cells.Item(327, 4).Value = 44446
$ws.Cells.Item(327, 12).Value = "Primera"
$ws.Cells.Item(327, 13).Value = 250
$ws.Cells.Item(327, 14).Value = 30000
$ws.Cells.Item(327, 15).Value = 30000
$ws.Cells.Item(327, 16).Value = 30000
$ws.Cells.Item(327, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(327, 18).Value = "Perú"
$ws.Cells.Item(327, 19).Value = 3000
$ws.Cells.Item(327, 20).Value = 10
$ws.Cells.Item(328, 4).Value = 44421
$ws.Cells.Item(328, 12).Value = "Especial"
$ws.Cells.Item(328, 13).Value = 150
$ws.Cells.Item(328, 14).Value = 35000
$ws.Cells.Item(328, 15).Value = 35000
$ws.Cells.Item(328, 16).Value = 35000
$ws.Cells.Item(328, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(328, 18).Value = "Perú"
$ws.Cells.Item(328, 19).Value = 3500
$ws.Cells.Item(328, 20).Value = 10
$ws.Cells.Item(329, 4).Value = 44421
$ws.Cells.Item(329, 12).Value = "Primera"
$ws.Cells.Item(329, 13).Value = 150
$ws.Cells.Item(329, 14).Value = 28000
$ws.Cells.Item(329, 15).Value = 28000
$ws.Cells.Item(329, 16).Value = 28000
$ws.Cells.Item(329, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(329, 18).Value = "Perú"
$ws.Cells.Item(329, 19).Value = 2800
$ws.Cells.Item(329, 20).Value = 10
$ws.Cells.Item(330, 4).Value = 44208
$ws.Cells.Item(330, 12).Value = "Primera"
$ws.Cells.Item(330, 13).Value = 200
$ws.Cells.Item(330, 14).Value = 5700
$ws.Cells.Item(330, 15).Value = 5800
$ws.Cells.Item(330, 16).Value = 5750
$ws.Cells.Item(330, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(330, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(330, 19).Value = 5750
$ws.Cells.Item(330, 20).Value = 1
$ws.Cells.Item(331, 4).Value = 44208
$ws.Cells.Item(331, 12).Value = "Segunda"
$ws.Cells.Item(331, 13).Value = 100
$ws.Cells.Item(331, 14).Value = 4800
$ws.Cells.Item(331, 15).Value = 4800
$ws.Cells.Item(331, 16).Value = 4800
$ws.Cells.Item(331, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(331, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(331, 19).Value = 4800
$ws.Cells.Item(331, 20).Value = 1
$ws.Cells.Item(332, 4).Value = 44355
$ws.Cells.Item(332, 12).Value = "Primera"
$ws.Cells.Item(332, 13).Value = 200
$ws.Cells.Item(332, 14).Value = 42000
$ws.Cells.Item(332, 15).Value = 43000
$ws.Cells.Item(332, 16).Value = 42500
$ws.Cells.Item(332, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(332, 18).Value = "Perú"
$ws.Cells.Item(332, 19).Value = 4250
$ws.Cells.Item(332, 20).Value = 10
$ws.Cells.Item(333, 4).Value = 44355
$ws.Cells.Item(333, 12).Value = "Primera"
$ws.Cells.Item(333, 13).Value = 200
$ws.Cells.Item(333, 14).Value = 7400
$ws.Cells.Item(333, 15).Value = 7500
$ws.Cells.Item(333, 16).Value = 7450
$ws.Cells.Item(333, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(333, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(333, 19).Value = 7450
$ws.Cells.Item(333, 20).Value = 1
$ws.Cells.Item(334, 4).Value = 44530
$ws.Cells.Item(334, 12).Value = "Primera"
$ws.Cells.Item(334, 13).Value = 400
$ws.Cells.Item(334, 14).Value = 4000
$ws.Cells.Item(334, 15).Value = 4100
$ws.Cells.Item(334, 16).Value = 4050
$ws.Cells.Item(334, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(334, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(334, 19).Value = 4050
$ws.Cells.Item(334, 20).Value = 1
$ws.Cells.Item(335, 4).Value = 44530
$ws.Cells.Item(335, 12).Value = "Segunda"
$ws.Cells.Item(335, 13).Value = 200
$ws.Cells.Item(335, 14).Value = 3500
$ws.Cells.Item(335, 15).Value = 3500
$ws.Cells.Item(335, 16).Value = 3500
$ws.Cells.Item(335, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(335, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(335, 19).Value = 3500
$ws.Cells.Item(335, 20).Value = 1
$ws.Cells.Item(336, 4).Value = 44483
$ws.Cells.Item(336, 12).Value = "1a nueva(o)"
$ws.Cells.Item(336, 13).Value = 200
$ws.Cells.Item(336, 14).Value = 4000
$ws.Cells.Item(336, 15).Value = 4200
$ws.Cells.Item(336, 16).Value = 4100
$ws.Cells.Item(336, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(336, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(336, 19).Value = 4100
$ws.Cells.Item(336, 20).Value = 1
$ws.Cells.Item(337, 4).Value = 44294
$ws.Cells.Item(337, 12).Value = "Primera"
$ws.Cells.Item(337, 13).Value = 60
$ws.Cells.Item(337, 14).Value = 50000
$ws.Cells.Item(337, 15).Value = 50000
$ws.Cells.Item(337, 16).Value = 50000
$ws.Cells.Item(337, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(337, 18).Value = "Perú"
$ws.Cells.Item(337, 19).Value = 5000
$ws.Cells.Item(337, 20).Value = 10
$ws.Cells.Item(338, 4).Value = 44294
$ws.Cells.Item(338, 12).Value = "Primera"
$ws.Cells.Item(338, 13).Value = 200
$ws.Cells.Item(338, 14).Value = 6300
$ws.Cells.Item(338, 15).Value = 6400
$ws.Cells.Item(338, 16).Value = 6350
$ws.Cells.Item(338, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(338, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(338, 19).Value = 6350
$ws.Cells.Item(338, 20).Value = 1
$ws.Cells.Item(339, 4).Value = 44617
$ws.Cells.Item(339, 12).Value = "Primera"
$ws.Cells.Item(339, 13).Value = 300
$ws.Cells.Item(339, 14).Value = 3900
$ws.Cells.Item(339, 15).Value = 4000
$ws.Cells.Item(339, 16).Value = 3950
$ws.Cells.Item(339, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(339, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(339, 19).Value = 3950
$ws.Cells.Item(339, 20).Value = 1
$ws.Cells.Item(340, 4).Value = 44617
$ws.Cells.Item(340, 12).Value = "Segunda"
$ws.Cells.Item(340, 13).Value = 150
$ws.Cells.Item(340, 14).Value = 3500
$ws.Cells.Item(340, 15).Value = 3500
$ws.Cells.Item(340, 16).Value = 3500
$ws.Cells.Item(340, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(340, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(340, 19).Value = 3500
$ws.Cells.Item(340, 20).Value = 1
$ws.Cells.Item(341, 4).Value = 44557
$ws.Cells.Item(341, 12).Value = "Primera"
$ws.Cells.Item(341, 13).Value = 400
$ws.Cells.Item(341, 14).Value = 3900
$ws.Cells.Item(341, 15).Value = 4000
$ws.Cells.Item(341, 16).Value = 3950
$ws.Cells.Item(341, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(341, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(341, 19).Value = 3950
$ws.Cells.Item(341, 20).Value = 1
$ws.Cells.Item(342, 4).Value = 44557
$ws.Cells.Item(342, 12).Value = "Segunda"
$ws.Cells.Item(342, 13).Value = 100
$ws.Cells.Item(342, 14).Value = 3500
$ws.Cells.Item(342, 15).Value = 3500
$ws.Cells.Item(342, 16).Value = 3500
$ws.Cells.Item(342, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(342, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(342, 19).Value = 3500
$ws.Cells.Item(342, 20).Value = 1
$ws.Cells.Item(343, 4).Value = 44264
$ws.Cells.Item(343, 12).Value = "Primera"
$ws.Cells.Item(343, 13).Value = 200
$ws.Cells.Item(343, 14).Value = 5800
$ws.Cells.Item(343, 15).Value = 6000
$ws.Cells.Item(343, 16).Value = 5900
$ws.Cells.Item(343, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(343, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(343, 19).Value = 5900
$ws.Cells.Item(343, 20).Value = 1
$ws.Cells.Item(344, 4).Value = 44264
$ws.Cells.Item(344, 12).Value = "Segunda"
$ws.Cells.Item(344, 13).Value = 100
$ws.Cells.Item(344, 14).Value = 4900
$ws.Cells.Item(344, 15).Value = 4900
$ws.Cells.Item(344, 16).Value = 4900
$ws.Cells.Item(344, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(344, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(344, 19).Value = 4900
$ws.Cells.Item(344, 20).Value = 1
$ws.Cells.Item(345, 4).Value = 44232
$ws.Cells.Item(345, 12).Value = "Primera"
$ws.Cells.Item(345, 13).Value = 250
$ws.Cells.Item(345, 14).Value = 5300
$ws.Cells.Item(345, 15).Value = 5300
$ws.Cells.Item(345, 16).Value = 5300
$ws.Cells.Item(345, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(345, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(345, 19).Value = 5300
$ws.Cells.Item(345, 20).Value = 1
$ws.Cells.Item(346, 4).Value = 44279
$ws.Cells.Item(346, 12).Value = "Primera"
$ws.Cells.Item(346, 13).Value = 80
$ws.Cells.Item(346, 14).Value = 6300
$ws.Cells.Item(346, 15).Value = 6400
$ws.Cells.Item(346, 16).Value = 6350
$ws.Cells.Item(346, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(346, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(346, 19).Value = 6350
$ws.Cells.Item(346, 20).Value = 1
$ws.Cells.Item(347, 4).Value = 44330
$ws.Cells.Item(347, 12).Value = "Primera"
$ws.Cells.Item(347, 13).Value = 200
$ws.Cells.Item(347, 14).Value = 46000
$ws.Cells.Item(347, 15).Value = 47000
$ws.Cells.Item(347, 16).Value = 46500
$ws.Cells.Item(347, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(347, 18).Value = "Perú"
$ws.Cells.Item(347, 19).Value = 4650
$ws.Cells.Item(347, 20).Value = 10
$ws.Cells.Item(348, 4).Value = 44330
$ws.Cells.Item(348, 12).Value = "Primera"
$ws.Cells.Item(348, 13).Value = 200
$ws.Cells.Item(348, 14).Value = 7000
$ws.Cells.Item(348, 15).Value = 7200
$ws.Cells.Item(348, 16).Value = 7100
$ws.Cells.Item(348, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(348, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(348, 19).Value = 7100
$ws.Cells.Item(348, 20).Value = 1
$ws.Cells.Item(349, 4).Value = 44572
$ws.Cells.Item(349, 12).Value = "Primera"
$ws.Cells.Item(349, 13).Value = 400
$ws.Cells.Item(349, 14).Value = 4000
$ws.Cells.Item(349, 15).Value = 4100
$ws.Cells.Item(349, 16).Value = 4050
$ws.Cells.Item(349, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(349, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(349, 19).Value = 4050
$ws.Cells.Item(349, 20).Value = 1
$ws.Cells.Item(350, 4).Value = 44572
$ws.Cells.Item(350, 12).Value = "Segunda"
$ws.Cells.Item(350, 13).Value = 200
$ws.Cells.Item(350, 14).Value = 3600
$ws.Cells.Item(350, 15).Value = 3600
$ws.Cells.Item(350, 16).Value = 3600
$ws.Cells.Item(350, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(350, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(350, 19).Value = 3600
$ws.Cells.Item(350, 20).Value = 1
$ws.Cells.Item(351, 4).Value = 44257
$ws.Cells.Item(351, 12).Value = "Primera"
$ws.Cells.Item(351, 13).Value = 200
$ws.Cells.Item(351, 14).Value = 5500
$ws.Cells.Item(351, 15).Value = 5700
$ws.Cells.Item(351, 16).Value = 5600
$ws.Cells.Item(351, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(351, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(351, 19).Value = 5600
$ws.Cells.Item(351, 20).Value = 1
$ws.Cells.Item(352, 4).Value = 44257
$ws.Cells.Item(352, 12).Value = "Segunda"
$ws.Cells.Item(352, 13).Value = 100
$ws.Cells.Item(352, 14).Value = 4300
$ws.Cells.Item(352, 15).Value = 4300
$ws.Cells.Item(352, 16).Value = 4300
$ws.Cells.Item(352, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(352, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(352, 19).Value = 4300
$ws.Cells.Item(352, 20).Value = 1
$ws.Cells.Item(353, 4).Value = 44236
$ws.Cells.Item(353, 12).Value = "Primera"
$ws.Cells.Item(353, 13).Value = 300
$ws.Cells.Item(353, 14).Value = 5400
$ws.Cells.Item(353, 15).Value = 5500
$ws.Cells.Item(353, 16).Value = 5450
$ws.Cells.Item(353, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(353, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(353, 19).Value = 5450
$ws.Cells.Item(353, 20).Value = 1
$ws.Cells.Item(354, 4).Value = 44236
$ws.Cells.Item(354, 12).Value = "Segunda"
$ws.Cells.Item(354, 13).Value = 150
$ws.Cells.Item(354, 14).Value = 4300
$ws.Cells.Item(354, 15).Value = 4300
$ws.Cells.Item(354, 16).Value = 4300
$ws.Cells.Item(354, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(354, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(354, 19).Value = 4300
$ws.Cells.Item(354, 20).Value = 1
$ws.Cells.Item(355, 4).Value = 44229
$ws.Cells.Item(355, 12).Value = "Primera"
$ws.Cells.Item(355, 13).Value = 200
$ws.Cells.Item(355, 14).Value = 5400
$ws.Cells.Item(355, 15).Value = 5500
$ws.Cells.Item(355, 16).Value = 5450
$ws.Cells.Item(355, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(355, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(355, 19).Value = 5450
$ws.Cells.Item(355, 20).Value = 1
$ws.Cells.Item(356, 4).Value = 44229
$ws.Cells.Item(356, 12).Value = "Segunda"
$ws.Cells.Item(356, 13).Value = 100
$ws.Cells.Item(356, 14).Value = 4200
$ws.Cells.Item(356, 15).Value = 4200
$ws.Cells.Item(356, 16).Value = 4200
$ws.Cells.Item(356, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(356, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(356, 19).Value = 4200
$ws.Cells.Item(356, 20).Value = 1
$ws.Cells.Item(357, 4).Value = 44299
$ws.Cells.Item(357, 12).Value = "Primera"
$ws.Cells.Item(357, 13).Value = 200
$ws.Cells.Item(357, 14).Value = 45000
$ws.Cells.Item(357, 15).Value = 46000
$ws.Cells.Item(357, 16).Value = 45500
$ws.Cells.Item(357, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(357, 18).Value = "Perú"
$ws.Cells.Item(357, 19).Value = 4550
$ws.Cells.Item(357, 20).Value = 10
$ws.Cells.Item(358, 4).Value = 44299
$ws.Cells.Item(358, 12).Value = "Primera"
$ws.Cells.Item(358, 13).Value = 200
$ws.Cells.Item(358, 14).Value = 6600
$ws.Cells.Item(358, 15).Value = 6700
$ws.Cells.Item(358, 16).Value = 6650
$ws.Cells.Item(358, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(358, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(358, 19).Value = 6650
$ws.Cells.Item(358, 20).Value = 1
$ws.Cells.Item(359, 4).Value = 44610
$ws.Cells.Item(359, 12).Value = "Primera"
$ws.Cells.Item(359, 13).Value = 300
$ws.Cells.Item(359, 14).Value = 4000
$ws.Cells.Item(359, 15).Value = 4000
$ws.Cells.Item(359, 16).Value = 4000
$ws.Cells.Item(359, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(359, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(359, 19).Value = 4000
$ws.Cells.Item(359, 20).Value = 1
$ws.Cells.Item(360, 4).Value = 44610
$ws.Cells.Item(360, 12).Value = "Segunda"
$ws.Cells.Item(360, 13).Value = 200
$ws.Cells.Item(360, 14).Value = 3500
$ws.Cells.Item(360, 15).Value = 3500
$ws.Cells.Item(360, 16).Value = 3500
$ws.Cells.Item(360, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(360, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(360, 19).Value = 3500
$ws.Cells.Item(360, 20).Value = 1
$ws.Cells.Item(361, 4).Value = 44399
$ws.Cells.Item(361, 12).Value = "Primera"
$ws.Cells.Item(361, 13).Value = 100
$ws.Cells.Item(361, 14).Value = 36000
$ws.Cells.Item(361, 15).Value = 36000
$ws.Cells.Item(361, 16).Value = 36000
$ws.Cells.Item(361, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(361, 18).Value = "Perú"
$ws.Cells.Item(361, 19).Value = 3600
$ws.Cells.Item(361, 20).Value = 10
$ws.Cells.Item(362, 4).Value = 44615
$ws.Cells.Item(362, 12).Value = "Primera"
$ws.Cells.Item(362, 13).Value = 80
$ws.Cells.Item(362, 14).Value = 3900
$ws.Cells.Item(362, 15).Value = 4000
$ws.Cells.Item(362, 16).Value = 3950
$ws.Cells.Item(362, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(362, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(362, 19).Value = 3950
$ws.Cells.Item(362, 20).Value = 1
$ws.Cells.Item(363, 4).Value = 44615
$ws.Cells.Item(363, 12).Value = "Segunda"
$ws.Cells.Item(363, 13).Value = 40
$ws.Cells.Item(363, 14).Value = 3500
$ws.Cells.Item(363, 15).Value = 3500
$ws.Cells.Item(363, 16).Value = 3500
$ws.Cells.Item(363, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(363, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(363, 19).Value = 3500
$ws.Cells.Item(363, 20).Value = 1
$ws.Cells.Item(364, 4).Value = 44522
$ws.Cells.Item(364, 12).Value = "Primera"
$ws.Cells.Item(364, 13).Value = 200
$ws.Cells.Item(364, 14).Value = 4000
$ws.Cells.Item(364, 15).Value = 4100
$ws.Cells.Item(364, 16).Value = 4050
$ws.Cells.Item(364, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(364, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(364, 19).Value = 4050
$ws.Cells.Item(364, 20).Value = 1
$ws.Cells.Item(365, 4).Value = 44522
$ws.Cells.Item(365, 12).Value = "Segunda"
$ws.Cells.Item(365, 13).Value = 100
$ws.Cells.Item(365, 14).Value = 3500
$ws.Cells.Item(365, 15).Value = 3500
$ws.Cells.Item(365, 16).Value = 3500
$ws.Cells.Item(365, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(365, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(365, 19).Value = 3500
$ws.Cells.Item(365, 20).Value = 1
$ws.Cells.Item(366, 4).Value = 44543
$ws.Cells.Item(366, 12).Value = "Primera"
$ws.Cells.Item(366, 13).Value = 200
$ws.Cells.Item(366, 14).Value = 4000
$ws.Cells.Item(366, 15).Value = 4100
$ws.Cells.Item(366, 16).Value = 4050
$ws.Cells.Item(366, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(366, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(366, 19).Value = 4050
$ws.Cells.Item(366, 20).Value = 1
$ws.Cells.Item(367, 4).Value = 44543
$ws.Cells.Item(367, 12).Value = "Segunda"
$ws.Cells.Item(367, 13).Value = 100
$ws.Cells.Item(367, 14).Value = 3600
$ws.Cells.Item(367, 15).Value = 3600
$ws.Cells.Item(367, 16).Value = 3600
$ws.Cells.Item(367, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(367, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(367, 19).Value = 3600
$ws.Cells.Item(367, 20).Value = 1
$ws.Cells.Item(368, 4).Value = 44390
$ws.Cells.Item(368, 12).Value = "2a nueva(o)"
$ws.Cells.Item(368, 13).Value = 350
$ws.Cells.Item(368, 14).Value = 4000
$ws.Cells.Item(368, 15).Value = 4000
$ws.Cells.Item(368, 16).Value = 4000
$ws.Cells.Item(368, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(368, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(368, 19).Value = 4000
$ws.Cells.Item(368, 20).Value = 1
$ws.Cells.Item(369, 4).Value = 44285
$ws.Cells.Item(369, 12).Value = "Primera"
$ws.Cells.Item(369, 13).Value = 200
$ws.Cells.Item(369, 14).Value = 6300
$ws.Cells.Item(369, 15).Value = 6400
$ws.Cells.Item(369, 16).Value = 6350
$ws.Cells.Item(369, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(369, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(369, 19).Value = 6350
$ws.Cells.Item(369, 20).Value = 1
$ws.Cells.Item(370, 4).Value = 44285
$ws.Cells.Item(370, 12).Value = "Segunda"
$ws.Cells.Item(370, 13).Value = 100
$ws.Cells.Item(370, 14).Value = 4990
$ws.Cells.Item(370, 15).Value = 4990
$ws.Cells.Item(370, 16).Value = 4990
$ws.Cells.Item(370, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(370, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(370, 19).Value = 4990
$ws.Cells.Item(370, 20).Value = 1
$ws.Cells.Item(371, 4).Value = 44498
$ws.Cells.Item(371, 12).Value = "1a nueva(o)"
$ws.Cells.Item(371, 13).Value = 400
$ws.Cells.Item(371, 14).Value = 4000
$ws.Cells.Item(371, 15).Value = 4200
$ws.Cells.Item(371, 16).Value = 4100
$ws.Cells.Item(371, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(371, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(371, 19).Value = 4100
$ws.Cells.Item(371, 20).Value = 1
$ws.Cells.Item(372, 4).Value = 44498
$ws.Cells.Item(372, 12).Value = "Segunda"
$ws.Cells.Item(372, 13).Value = 200
$ws.Cells.Item(372, 14).Value = 3600
$ws.Cells.Item(372, 15).Value = 3600
$ws.Cells.Item(372, 16).Value = 3600
$ws.Cells.Item(372, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(372, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(372, 19).Value = 3600
$ws.Cells.Item(372, 20).Value = 1
$ws.Cells.Item(373, 4).Value = 44418
$ws.Cells.Item(373, 12).Value = "Especial"
$ws.Cells.Item(373, 13).Value = 200
$ws.Cells.Item(373, 14).Value = 35000
$ws.Cells.Item(373, 15).Value = 35000
$ws.Cells.Item(373, 16).Value = 35000
$ws.Cells.Item(373, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(373, 18).Value = "Perú"
$ws.Cells.Item(373, 19).Value = 3500
$ws.Cells.Item(373, 20).Value = 10
$ws.Cells.Item(374, 4).Value = 44418
$ws.Cells.Item(374, 12).Value = "Primera"
$ws.Cells.Item(374, 13).Value = 150
$ws.Cells.Item(374, 14).Value = 28000
$ws.Cells.Item(374, 15).Value = 28000
$ws.Cells.Item(374, 16).Value = 28000
$ws.Cells.Item(374, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(374, 18).Value = "Perú"
$ws.Cells.Item(374, 19).Value = 2800
$ws.Cells.Item(374, 20).Value = 10
$ws.Cells.Item(375, 4).Value = 44595
$ws.Cells.Item(375, 12).Value = "Primera"
$ws.Cells.Item(375, 13).Value = 200
$ws.Cells.Item(375, 14).Value = 4000
$ws.Cells.Item(375, 15).Value = 4100
$ws.Cells.Item(375, 16).Value = 4050
$ws.Cells.Item(375, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(375, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(375, 19).Value = 4050
$ws.Cells.Item(375, 20).Value = 1
$ws.Cells.Item(376, 4).Value = 44595
$ws.Cells.Item(376, 12).Value = "Segunda"
$ws.Cells.Item(376, 13).Value = 100
$ws.Cells.Item(376, 14).Value = 3500
$ws.Cells.Item(376, 15).Value = 3500
$ws.Cells.Item(376, 16).Value = 3500
$ws.Cells.Item(376, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(376, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(376, 19).Value = 3500
$ws.Cells.Item(376, 20).Value = 1
$ws.Cells.Item(377, 4).Value = 44628
$ws.Cells.Item(377, 12).Value = "Primera"
$ws.Cells.Item(377, 13).Value = 600
$ws.Cells.Item(377, 14).Value = 4300
$ws.Cells.Item(377, 15).Value = 4400
$ws.Cells.Item(377, 16).Value = 4350
$ws.Cells.Item(377, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(377, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(377, 19).Value = 4350
$ws.Cells.Item(377, 20).Value = 1
$ws.Cells.Item(378, 4).Value = 44335
$ws.Cells.Item(378, 12).Value = "Primera"
$ws.Cells.Item(378, 13).Value = 120
$ws.Cells.Item(378, 14).Value = 7000
$ws.Cells.Item(378, 15).Value = 7200
$ws.Cells.Item(378, 16).Value = 7100
$ws.Cells.Item(378, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(378, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(378, 19).Value = 7100
$ws.Cells.Item(378, 20).Value = 1
$ws.Cells.Item(379, 4).Value = 44552
$ws.Cells.Item(379, 12).Value = "Primera"
$ws.Cells.Item(379, 13).Value = 200
$ws.Cells.Item(379, 14).Value = 3900
$ws.Cells.Item(379, 15).Value = 4000
$ws.Cells.Item(379, 16).Value = 3950
$ws.Cells.Item(379, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(379, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(379, 19).Value = 3950
$ws.Cells.Item(379, 20).Value = 1
$ws.Cells.Item(380, 4).Value = 44552
$ws.Cells.Item(380, 12).Value = "Segunda"
$ws.Cells.Item(380, 13).Value = 100
$ws.Cells.Item(380, 14).Value = 3500
$ws.Cells.Item(380, 15).Value = 3500
$ws.Cells.Item(380, 16).Value = 3500
$ws.Cells.Item(380, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(380, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(380, 19).Value = 3500
$ws.Cells.Item(380, 20).Value = 1
$ws.Cells.Item(381, 1).Value = 4
$ws.Cells.Item(381, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(381, 3).Value = "Los Lagos"
$ws.Cells.Item(381, 4).Value = 44544
$ws.Cells.Item(381, 5).Value = 10
$ws.Cells.Item(381, 6).Value = "Fruta"
$ws.Cells.Item(381, 7).Value = 100106
$ws.Cells.Item(381, 8).Value = "Oleaginosos"
$ws.Cells.Item(381, 9).Value = 100106002
$ws.Cells.Item(381, 10).Value = "Palta"
$ws.Cells.Item(381, 11).Value = "Hass"
$ws.Cells.Item(381, 12).Value = "Primera"
$ws.Cells.Item(381, 13).Value = 400
$ws.Cells.Item(381, 14).Value = 4000
$ws.Cells.Item(381, 15).Value = 4100
$ws.Cells.Item(381, 16).Value = 4050
$ws.Cells.Item(381, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(381, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(381, 19).Value = 4050
$ws.Cells.Item(381, 20).Value = 1
$ws.Cells.Item(381, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(382, 1).Value = 4
$ws.Cells.Item(382, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(382, 3).Value = "Los Lagos"
$ws.Cells.Item(382, 4).Value = 44544
$ws.Cells.Item(382, 5).Value = 10
$ws.Cells.Item(382, 6).Value = "Fruta"
$ws.Cells.Item(382, 7).Value = 100106
$ws.Cells.Item(382, 8).Value = "Oleaginosos"
$ws.Cells.Item(382, 9).Value = 100106002
$ws.Cells.Item(382, 10).Value = "Palta"
$ws.Cells.Item(382, 11).Value = "Hass"
$ws.Cells.Item(382, 12).Value = "Segunda"
$ws.Cells.Item(382, 13).Value = 200
$ws.Cells.Item(382, 14).Value = 3500
$ws.Cells.Item(382, 15).Value = 3500
$ws.Cells.Item(382, 16).Value = 3500
$ws.Cells.Item(382, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(382, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(382, 19).Value = 3500
$ws.Cells.Item(382, 20).Value = 1
$ws.Cells.Item(382, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

Write-Output "done"
